$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing "%" from the Accuracy column (column B, rows 2-11),
# keeping the cell content as text (not converting to a numeric value).
$values = @{
    2  = "59.53"
    3  = "60.14"
    4  = "57.76"
    5  = "53.38"
    6  = "59.81"
    7  = "52.54"
    8  = "60.50"
    9  = "59.80"
    10 = "53.40"
    11 = "59.97"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 2)
    # Force text format so Excel doesn't reinterpret the numeric-looking
    # string as a number, then restore the original (default) style so
    # no extra formatting is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$row]
    $cell.Style = "Normal"
}
